$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.331999999999999
$ws.Range("D21").Value = -7.805
$ws.Range("D23").Value = -7.455
$ws.Range("D25").Value = -8.388999999999999
